# Updates the cryptos price/volume table with the latest scraped values.
# Cells whose new text is a "plain" decimal number (e.g. "32.58") are written
# with a leading apostrophe so Excel keeps them as literal text (matching the
# source data's inlineStr cells) instead of silently converting them to
# floating point numbers; the style is then reset to "Normal" so no stray
# number-format/quote-prefix style sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.346.15'
$ws.Range('E2').Value = '  +0.65%  '
$ws.Range('D3').Value = '1.786.27'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''225.99'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').Value = '''0.553'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.13%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '''32.58'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.49%  '
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('E10').Value = '  +0.42%  '
$ws.Range('D11').Value = '''0.0946'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.45%  '
$ws.Range('D12').Value = '2.044.92'
$ws.Range('D13').Value = '1.786.94'
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('E14').Value = '  +0.50%  '
$ws.Range('D15').Value = '''0.631'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.60%  '
$ws.Range('D16').Value = '34.355.00'
$ws.Range('E16').Value = '  +0.71%  '
$ws.Range('D17').Value = '''4.27'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.30%  '
$ws.Range('D18').Value = '''68.17'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.89%  '
$ws.Range('E19').Value = '  +0.94%  '
$ws.Range('D20').Value = '''243.89'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.49%  '
$ws.Range('D21').Value = '''11.15'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.50%  '
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').Value = '''4.13'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.61%  '
$ws.Range('E24').Value = '  +1.91%  '
$ws.Range('D25').Value = '''165.81'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.22%  '
$ws.Range('D26').Value = '''7.27'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.09%  '
$ws.Range('D27').Value = '''16.46'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.13%  '
$ws.Range('E28').Value = '  +0.97%  '
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('D30').Value = '''3.96'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +6.78%  '
$ws.Range('E31').Value = '  +1.35%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '''1.23'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.64%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''3.79'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.40%  '
$ws.Range('D34').Value = '''1.80'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('E35').Value = '  +5.53%  '
$ws.Range('D36').Value = '1.404.43'
$ws.Range('E36').Value = '  -2.55%  '
$ws.Range('D37').Value = '''0.676'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +3.65%  '
$ws.Range('E38').Value = '  +2.40%  '
$ws.Range('E39').Value = '  -0.54%  '
$ws.Range('D40').Value = '''84.50'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +3.68%  '
$ws.Range('E41').Value = '  +0.77%  '
$ws.Range('E42').Value = '  +2.52%  '
$ws.Range('E43').Value = '  +2.59%  '
$ws.Range('D44').Value = '''13.80'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.15%  '
$ws.Range('D45').Value = '''0.0526'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.48%  '
$ws.Range('E46').Value = '  +3.08%  '
$ws.Range('E47').Value = '  +0.16%  '
$ws.Range('D48').Value = '1.945.31'
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('D49').Value = '''105.04'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('E51').Value = '  -2.55%  '
